# Update titles for closed bugs in the "closed bugs in last iteration" sheet.
# Only the Title (column B) text changes for a handful of rows; IDs (column A)
# and State (column C) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("closed bugs in last iteration")

$ws.Range("B2").Value = '"Assigned to" user selector in Test Plans sometimes resets itself back to current user after being cleared'
$ws.Range("B3").Value = "Remove inferred bus type value for an asset when bus type is not defined in the client "
$ws.Range("B5").Value = "Asset Location tab should update query request when changing the current asset"
$ws.Range("B14").Value = "FeedService.WebAPI - CVE-2021-26701"
$ws.Range("B22").Value = "Alarm and DFS example continuation tokens in Swagger UI are invalid"
$ws.Range("B24").Value = "RBAC integrator container doesn't emit the duration metric when the job fails"
$ws.Range("B43").Value = "Users can't download files larger than 10MB from the files grid"
